$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update D (Price) and E (Volume(1h)) columns for rows 2-47 with refreshed values.
# Cells whose new Price text is purely numeric-looking are forced to Text format
# first so Excel stores them as strings (matching the source data) rather than numbers.

$ws.Cells.Item(2, 4).Value = "26.842.44"
$ws.Cells.Item(2, 5).Value = "  +0.12%  "
$ws.Cells.Item(3, 4).Value = "1.643.83"
$ws.Cells.Item(3, 5).Value = "  -0.11%  "
$ws.Cells.Item(4, 5).Value = "  -0.38%  "
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "218.57"
$ws.Cells.Item(5, 5).Value = "  +0.86%  "
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "0.502"
$ws.Cells.Item(6, 5).Value = "  +0.39%  "
$ws.Cells.Item(7, 5).Value = "  -0.45%  "
$ws.Cells.Item(8, 5).Value = "  -0.19%  "
$ws.Cells.Item(9, 5).Value = "  -0.91%  "
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "19.28"
$ws.Cells.Item(10, 5).Value = "  +0.47%  "
$ws.Cells.Item(11, 5).Value = "  +0.82%  "
$ws.Cells.Item(12, 4).Value = "1.872.20"
$ws.Cells.Item(12, 5).Value = "  -0.12%  "
$ws.Cells.Item(13, 4).Value = "1.635.84"
$ws.Cells.Item(13, 5).Value = "  -0.60%  "
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "4.16"
$ws.Cells.Item(14, 5).Value = "  -0.22%  "
$ws.Cells.Item(15, 5).Value = "  -0.47%  "
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = "65.37"
$ws.Cells.Item(16, 5).Value = "  +1.48%  "
$ws.Cells.Item(17, 4).Value = "26.840.66"
$ws.Cells.Item(17, 5).Value = "  +0.09%  "
$ws.Cells.Item(18, 4).Value = "0.0₃0736"
$ws.Cells.Item(18, 5).Value = "  -0.36%  "
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "216.13"
$ws.Cells.Item(19, 5).Value = "  +1.18%  "
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "6.63"
$ws.Cells.Item(21, 5).Value = "  +5.32%  "
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "4.38"
$ws.Cells.Item(22, 5).Value = "  +0.22%  "
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "2.36"
$ws.Cells.Item(23, 5).Value = "  -0.99%  "
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "9.22"
$ws.Cells.Item(24, 5).Value = "  -1.61%  "
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "147.72"
$ws.Cells.Item(25, 5).Value = "  +1.84%  "
$ws.Cells.Item(26, 5).Value = "  -0.51%  "
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "0.119"
$ws.Cells.Item(27, 5).Value = "  +0.00%  "
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "7.14"
$ws.Cells.Item(28, 5).Value = "  +0.83%  "
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "15.77"
$ws.Cells.Item(29, 5).Value = "  +0.62%  "
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "0.0510"
$ws.Cells.Item(30, 5).Value = "  -0.04%  "
$ws.Cells.Item(31, 5).Value = "  +0.98%  "
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "3.39"
$ws.Cells.Item(32, 5).Value = "  +2.15%  "
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "2.99"
$ws.Cells.Item(33, 5).Value = "  -0.60%  "
$ws.Cells.Item(34, 5).Value = "  +0.97%  "
$ws.Cells.Item(35, 4).Value = "1.269.56"
$ws.Cells.Item(35, 5).Value = "  -1.41%  "
$ws.Cells.Item(36, 5).Value = "  +0.27%  "
$ws.Cells.Item(37, 5).Value = "  +0.96%  "
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "0.533"
$ws.Cells.Item(38, 5).Value = "  -1.05%  "
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "0.818"
$ws.Cells.Item(39, 5).Value = "  -0.99%  "
$ws.Cells.Item(40, 5).Value = "  -0.41%  "
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "0.805"
$ws.Cells.Item(41, 5).Value = "  -0.66%  "
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "5.35"
$ws.Cells.Item(42, 5).Value = "  -0.18%  "
$ws.Cells.Item(43, 4).Value = "1.782.10"
$ws.Cells.Item(43, 5).Value = "  -0.80%  "
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "2.14"
$ws.Cells.Item(44, 5).Value = "  -4.49%  "
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "92.79"
$ws.Cells.Item(45, 5).Value = "  +1.44%  "
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "61.33"
$ws.Cells.Item(46, 5).Value = "  +1.29%  "
$ws.Cells.Item(47, 5).Value = "  +0.26%  "

# Rows 48-51: BabyDogeCoin inserted, shifting Cronos/Algorand/EnergySwap down; USDD drops off the list
$ws.Cells.Item(48, 2).Value = "BabyDogeCoin"
$ws.Cells.Item(48, 3).Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Cells.Item(48, 4).Value = "0.0₆0102"
$ws.Cells.Item(48, 5).Value = "  -1.34%  "

$ws.Cells.Item(49, 2).Value = "Cronos"
$ws.Cells.Item(49, 3).Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "0.0516"
$ws.Cells.Item(49, 5).Value = "  -0.43%  "

$ws.Cells.Item(50, 2).Value = "Algorand"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "0.0967"
$ws.Cells.Item(50, 5).Value = "  -1.15%  "

$ws.Cells.Item(51, 2).Value = "EnergySwap"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "7.53"
$ws.Cells.Item(51, 5).Value = "  -1.76%  "
